$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend formatting (style) for the new rows 17-34 by copying the format
# of the last existing data row (row 16) down across the new range.
$ws.Range("A16:C16").Copy()
$ws.Range("A17:C34").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Write out the full, reordered/updated abbreviation table (rows 1-34).
$ws.Cells.Item(1, 1).Value = "AkNr"
$ws.Cells.Item(1, 2).Value = "KlartextDe"
$ws.Cells.Item(1, 3).Value = "KlartextEn"

$ws.Cells.Item(2, 1).Value = "SDG"
$ws.Cells.Item(2, 2).Value = "Ziele für Nachhaltige Entwicklung (Sustainable Development Goals)"
$ws.Cells.Item(2, 3).Value = "Sustainable Development Goals"

$ws.Cells.Item(3, 1).Value = "OECD"
$ws.Cells.Item(3, 2).Value = "Organisation für wirtschaftliche Zusammenarbeit und Entwicklung (Organisation for Economic Co-operation and Development)"
$ws.Cells.Item(3, 3).Value = "Organisation for Economic Co-operation and Development"

$ws.Cells.Item(4, 1).Value = "ODA"
$ws.Cells.Item(4, 2).Value = "Öffentliche Entwicklungsausgaben (official development assistance)"
$ws.Cells.Item(4, 3).Value = "Official development assistance"

$ws.Cells.Item(5, 1).Value = "kg/ha"
$ws.Cells.Item(5, 2).Value = "Kilogramm pro Hektar"
$ws.Cells.Item(5, 3).Value = "Kilogram per hectare"

$ws.Cells.Item(6, 1).Value = "kg"
$ws.Cells.Item(6, 2).Value = "Kilogramm"
$ws.Cells.Item(6, 3).Value = "Kilogram"

$ws.Cells.Item(7, 1).Value = "ha"
$ws.Cells.Item(7, 2).Value = "Hektar"
$ws.Cells.Item(7, 3).Value = "Hectare"

$ws.Cells.Item(8, 1).Value = "FAO"
$ws.Cells.Item(8, 2).Value = "Ernährungs- und Landwirtschaftsorganisation der Vereinten Nationen (Food and Agriculture Organization)"
$ws.Cells.Item(8, 3).Value = "Food and Agriculture Organization"

$ws.Cells.Item(9, 1).Value = "EU-SILC"
$ws.Cells.Item(9, 2).Value = "Statistik über Einkommen und Lebensbedingungen (Statistics on Income and Living Conditions)"
$ws.Cells.Item(9, 3).Value = "Statistics on Income and Living Conditions"

$ws.Cells.Item(10, 1).Value = "EU-28"
$ws.Cells.Item(10, 2).Value = "Europäische Union mit 28 Mitgliedsstaaten"
$ws.Cells.Item(10, 3).Value = "European Union consisting of 28 member states"

$ws.Cells.Item(11, 1).Value = "EU-27"
$ws.Cells.Item(11, 2).Value = "Europäische Union mit 27 Mitgliedsstaaten"
$ws.Cells.Item(11, 3).Value = "European Union consisting of 27 member states"

$ws.Cells.Item(12, 1).Value = "EU"
$ws.Cells.Item(12, 2).Value = "Europäische Union"
$ws.Cells.Item(12, 3).Value = "European Union"

$ws.Cells.Item(13, 1).Value = "bzw."
$ws.Cells.Item(13, 2).Value = "beziehungsweise"
$ws.Cells.Item(13, 3).Value = ""

$ws.Cells.Item(14, 1).Value = "BMZ"
$ws.Cells.Item(14, 2).Value = "Bundesministerium für wirtschaftliche Zusammenarbeit und Entwicklung"
$ws.Cells.Item(14, 3).Value = "Federal Ministry for Economic Cooperation and Developmen"

$ws.Cells.Item(15, 1).Value = "BMEL"
$ws.Cells.Item(15, 2).Value = "Bundesministeriums für Ernährung und Landwirtschaft"
$ws.Cells.Item(15, 3).Value = "Federal Ministry of Food and Agriculture"

$ws.Cells.Item(16, 1).Value = "BLE"
$ws.Cells.Item(16, 2).Value = "Bundesanstalt für Landwirtschaft und Ernährung"
$ws.Cells.Item(16, 3).Value = "Federal Office for Agriculture and Food"

$ws.Cells.Item(17, 1).Value = "SDGs"
$ws.Cells.Item(17, 2).Value = "Ziele für Nachhaltige Entwicklung (Sustainable Development Goals)"
$ws.Cells.Item(17, 3).Value = "Sustainable Development Goals"

$ws.Cells.Item(18, 1).Value = "BMI"
$ws.Cells.Item(18, 2).Value = "Body Mass Index"
$ws.Cells.Item(18, 3).Value = "Body Mass Index"

$ws.Cells.Item(19, 1).Value = "kg/m²"
$ws.Cells.Item(19, 2).Value = "Kilogramm pro Quadratmeter"
$ws.Cells.Item(19, 3).Value = "Kilogram per square meter"

$ws.Cells.Item(20, 1).Value = "P90"
$ws.Cells.Item(20, 2).Value = "90. Perzentil"
$ws.Cells.Item(20, 3).Value = "90th percentile"

$ws.Cells.Item(21, 1).Value = "P97"
$ws.Cells.Item(21, 2).Value = "97. Perzentil"
$ws.Cells.Item(21, 3).Value = "97th percentile"

$ws.Cells.Item(22, 1).Value = "RKI"
$ws.Cells.Item(22, 2).Value = "Robert Koch-Institut"
$ws.Cells.Item(22, 3).Value = "Robert Koch-Institute"

$ws.Cells.Item(23, 1).Value = "KiGGS"
$ws.Cells.Item(23, 2).Value = "Studie zur Gesundheit von Kindern und Jugendlichen in Deutschland"
$ws.Cells.Item(23, 3).Value = "Study on the health of children and adolescents in Germany"

$ws.Cells.Item(24, 1).Value = "SES"
$ws.Cells.Item(24, 2).Value = "Sozioökonomischer Status"
$ws.Cells.Item(24, 3).Value = "Socioeconomic status"

$ws.Cells.Item(25, 1).Value = "WHO"
$ws.Cells.Item(25, 2).Value = "Weltgesundheitsorganisation (World Health Organization)"
$ws.Cells.Item(25, 3).Value = "World Health Organization"

$ws.Cells.Item(26, 1).Value = "SO₂"
$ws.Cells.Item(26, 2).Value = "Schwefeldioxid"
$ws.Cells.Item(26, 3).Value = "Sulphur dioxide"

$ws.Cells.Item(27, 1).Value = "NOₓ"
$ws.Cells.Item(27, 2).Value = "Stickstoffoxid"
$ws.Cells.Item(27, 3).Value = "Nitrogen oxides"

$ws.Cells.Item(28, 1).Value = "NH₃"
$ws.Cells.Item(28, 2).Value = "Ammoniak"
$ws.Cells.Item(28, 3).Value = "Ammonia"

$ws.Cells.Item(29, 1).Value = "NMVOC"
$ws.Cells.Item(29, 2).Value = "Flüchtige organische Verbindungen (non-methane volatile organic compounds)"
$ws.Cells.Item(29, 3).Value = "non-methane volatile organic compounds"

$ws.Cells.Item(30, 1).Value = "PM₂,₅"
$ws.Cells.Item(30, 2).Value = "Feinstaub"
$ws.Cells.Item(30, 3).Value = "Particulate matter"

$ws.Cells.Item(31, 1).Value = "CLRTAP"
$ws.Cells.Item(31, 2).Value = "Genfer Luftreinhaltekonvention (Convention on Long-Range Transboundary Air Pollution)"
$ws.Cells.Item(31, 3).Value = "Convention on Long-Range Transboundary Air Pollution"

$ws.Cells.Item(32, 1).Value = "NEC"
$ws.Cells.Item(32, 2).Value = "Richtlinie über nationale Emissionshöchstmengen für bestimmte Luftschadstoffe (National Emission Ceilings Directive)"
$ws.Cells.Item(32, 3).Value = "National Emission Ceilings Directive"

$ws.Cells.Item(33, 1).Value = "PM₂.₅"
$ws.Cells.Item(33, 2).Value = "Feinstaub"
$ws.Cells.Item(33, 3).Value = "Particulate matter"

$ws.Cells.Item(34, 1).Value = "NMVOCs"
$ws.Cells.Item(34, 2).Value = "Flüchtige organische Verbindungen (non-methane volatile organic compounds)"
$ws.Cells.Item(34, 3).Value = "Non-methane volatile organic compounds"
